$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new date column ("15-nov") is being appended right after the existing
# last column (CL, "14-nov"). Copy CL's formatting/values into the new CM
# column first (so styles match exactly, like Excel does when you extend a
# table by filling right), then overwrite the header and the handful of
# data cells that actually differ from the prior column's values.
$ws.Range("CL1:CL11").Copy($ws.Range("CM1:CM11"))

# New header label for the appended column.
$ws.Range("CM1").Value = "15-nov"

# Data values for the new column (rows 2-11, family Familia entries).
$ws.Range("CM2").Value = 9
$ws.Range("CM3").Value = 10
$ws.Range("CM4").Value = 7
$ws.Range("CM5").Value = 9
$ws.Range("CM6").Value = 12
$ws.Range("CM7").Value = 7
$ws.Range("CM8").Value = 14
$ws.Range("CM9").Value = 13
$ws.Range("CM10").Value = 21
$ws.Range("CM11").Value = 0

# Match the final selection state recorded in the workbook.
$ws.Range("CM11").Select()
